$d = $word.ActiveDocument

$pairs = @(
    @("586×7=4102", "615×5=3075"),
    @("420×8=3360", "128×2=256"),
    @("998×7=6986", "682×4=2728"),
    @("791×4=3164", "131×6=786"),
    @("269×3=807",  "774×3=2322"),
    @("127×5=635",  "700×5=3500"),
    @("636×6=3816", "670×4=2680"),
    @("783×7=5481", "819×7=5733"),
    @("432×9=3888", "900×8=7200"),
    @("563×3=1689", "289×9=2601"),
    @("665×6=3990", "720×8=5760"),
    @("891×8=7128", "362×8=2896"),
    @("447×2=894",  "698×6=4188"),
    @("346×6=2076", "970×2=1940"),
    @("759×8=6072", "923×5=4615"),
    @("908×7=6356", "594×5=2970"),
    @("898×7=6286", "621×2=1242"),
    @("675×9=6075", "914×2=1828"),
    @("923×2=1846", "631×3=1893"),
    @("555×9=4995", "490×8=3920"),
    @("720×7=5040", "959×8=7672"),
    @("478×5=2390", "404×7=2828"),
    @("840×6=5040", "306×8=2448"),
    @("575×3=1725", "764×5=3820"),
    @("619×9=5571", "179×7=1253")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
